# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "Strike#" column (G) is replaced with the newly-calculated "K" values.
# Only column G (the 7th column) changes; rows 2-69 on the single sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value (was Strike#)
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 2
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 0
    45 = 1
    46 = 2
    47 = 0
    48 = 2
    49 = 0
    50 = 1
    51 = 2
    52 = 1
    53 = 2
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    58 = 0
    59 = 2
    60 = 2
    61 = 0
    62 = 3
    63 = 1
    64 = 2
    65 = 1
    66 = 2
    67 = 2
    68 = 1
    69 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
